# Update to register answers
$wb = $excel.ActiveWorkbook

# --- "Registros" sheet: row 7 (Maria Aparecida Silva) ---
# The sheet used to hold 1/0 correction flags per question; it now holds the
# actually registered answer ("a"/"b"/"c"/"d" or "-" for blank) per question,
# plus an updated score (Nota) and an updated ethnicity value.
$ws = $wb.Worksheets.Item("Registros")

# B7 ("Nota") must stay text ("4"), not become a number - format the cell as
# text first so Excel keeps it as a string like all the other cells in the
# column.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "4"

$ws.Range("C7").Value = "no"
$ws.Range("D7").Value = "a"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "c"
$ws.Range("G7").Value = "c"
$ws.Range("H7").Value = "d"
$ws.Range("I7").Value = "b"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "c"
$ws.Range("L7").Value = "d"
$ws.Range("M7").Value = "b"
$ws.Range("N7").Value = "d"
$ws.Range("O7").Value = "d"
$ws.Range("P7").Value = "-"
$ws.Range("Q7").Value = "b"
$ws.Range("R7").Value = "c"
$ws.Range("S7").Value = "-"
$ws.Range("T7").Value = "d"
$ws.Range("U7").Value = "c"
$ws.Range("V7").Value = "b"
$ws.Range("W7").Value = "c"

# --- "Gabarito" sheet: used range shrinks from A1:F21 to A1:C21 ---
# (the sheet only ever used columns A-C; no cell data changes there, just
# re-activate/select it like the saved view shows)
$ws2 = $wb.Worksheets.Item("Gabarito")
$ws2.Activate()
$ws2.Range("F20").Select()
